$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: "A" across B2:G2 ---
$ws.Range("B2:G2").Value = "A"

# --- Serial numbers down column B ---
$ws.Range("B3").Value = 1
$ws.Range("B4").Formula = "=B3+1"
$ws.Range("B5:B18").Formula = "=B4+1"

# --- Borders around the header + number column (B2:B18, C2:G2) ---
$ws.Range("B2:B18").Borders.LineStyle = 1
$ws.Range("C2:G2").Borders.LineStyle = 1

# --- Body table: bordered + yellow filled block ---
$ws.Range("C3:G18").Borders.LineStyle = 1
$ws.Range("C3:G19").Interior.Color = 65535

# --- Selection, matching the saved view state ---
$ws.Range("J10").Select() | Out-Null

$wb.Save()
